# Update "Last Updated" timestamp on the Metadata sheet
$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "05 Nov 2025, 12:18 PM"

# Update the "1 Year" column (F) values on the Industry Analysis sheet
$ia = $wb.Worksheets.Item("Industry Analysis")
$ia.Range("F2").Value = 21.3
$ia.Range("F3").Value = -4.3927
$ia.Range("F4").Value = 35.9445
$ia.Range("F5").Value = -51.0482
$ia.Range("F6").Value = 57.2275
$ia.Range("F7").Value = -9.640700000000001
$ia.Range("F8").Value = -6.1449
$ia.Range("F9").Value = 36.9733
$ia.Range("F10").Value = -4.7026
$ia.Range("F11").Value = 46.5317
$ia.Range("F12").Value = -2.102
$ia.Range("F13").Value = 17.4681
$ia.Range("F14").Value = -33.0245
$ia.Range("F15").Value = 1.0205
$ia.Range("F16").Value = 2.0426
$ia.Range("F17").Value = -16.2411
$ia.Range("F18").Value = 7.4627
$ia.Range("F19").Value = -25.798
$ia.Range("F20").Value = 47.7485
$ia.Range("F21").Value = 19.5587
$ia.Range("F22").Value = 76.5603
$ia.Range("F23").Value = -54.2675
$ia.Range("F24").Value = -0.8811
$ia.Range("F25").Value = 4.8518
$ia.Range("F26").Value = 3.6831
$ia.Range("F27").Value = -34.0874
$ia.Range("F28").Value = -11.9893
$ia.Range("F29").Value = -12.994
$ia.Range("F30").Value = 25.5415
$ia.Range("F31").Value = 56.5088
$ia.Range("F32").Value = 2.0908
$ia.Range("F33").Value = -4.7193
$ia.Range("F34").Value = 22.8807
$ia.Range("F35").Value = 5.3359
$ia.Range("F36").Value = -5.1995
$ia.Range("F37").Value = -5.6238
$ia.Range("F38").Value = -22.595
$ia.Range("F39").Value = 10.8405
$ia.Range("F40").Value = -7.5963
$ia.Range("F41").Value = -4.552
$ia.Range("F42").Value = 22.3098
$ia.Range("F43").Value = 14.0694
$ia.Range("F44").Value = -9.6066
$ia.Range("F45").Value = 27.639
$ia.Range("F46").Value = -6.3484
$ia.Range("F47").Value = -40.5302
$ia.Range("F48").Value = -29.7988
$ia.Range("F49").Value = -24.0791
$ia.Range("F50").Value = -49.1803
$ia.Range("F51").Value = -51.6023
$ia.Range("F52").Value = -34.4756
$ia.Range("F53").Value = -11.5478
$ia.Range("F54").Value = -2.3796
$ia.Range("F55").Value = -15.4382
$ia.Range("F56").Value = -27.6987
$ia.Range("F57").Value = -27.1559
$ia.Range("F58").Value = -2.1585
$ia.Range("F59").Value = -23.0964
$ia.Range("F60").Value = -13.3217
$ia.Range("F61").Value = -8.1496
$ia.Range("F62").Value = -16.0695
$ia.Range("F63").Value = -12.5465
$ia.Range("F64").Value = 47.7264
$ia.Range("F65").Value = -42.4232
$ia.Range("F66").Value = 11.3291
$ia.Range("F67").Value = 14.3746
$ia.Range("F68").Value = 32.6702
$ia.Range("F69").Value = -17.0097
$ia.Range("F70").Value = -13.5162
$ia.Range("F71").Value = 11.4259
$ia.Range("F72").Value = 2.6754
$ia.Range("F73").Value = -11.1574
$ia.Range("F74").Value = -13.2502
$ia.Range("F75").Value = 24.7078
$ia.Range("F76").Value = 53.3554
